$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.381.05'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.37%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.867.32'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.64%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '339.03'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.08%  '
$ws.Range("E6").Value = '  -0.06%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4695'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.85%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3959'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.75%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '47.19'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08002'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.67%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.9995'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.45%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.86'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.47%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.864.47'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.09%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.991'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.72%  '
$ws.Range("E15").Value = '  +2.96%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.17'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.92%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.002'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.19%  '
$ws.Range("E18").Value = '  +1.21%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06613'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.36%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.56'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.71%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.002'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.02%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '28.392.60'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.41%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.452'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.13%  '
$ws.Range("E24").Value = '  +1.78%  '
$ws.Range("E25").Value = '  -1.32%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.090.71'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.27%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '160.52'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.14%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.78'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.33%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.127'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.00%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.489'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.10%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '120.02'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.23%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9664'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.34%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09485'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.99%  '
$ws.Range("E34").Value = '  +0.13%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.343'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.24%  '
$ws.Range("E36").Value = '  +3.96%  '
$ws.Range("E37").Value = '  +2.53%  '
$ws.Range("E38").Value = '  +2.65%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.348'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.43%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.188'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.32%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5932'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.18%  '
$ws.Range("E42").Value = '  +0.03%  '
$ws.Range("E43").Value = '  +1.53%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.32'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.46%  '
$ws.Range("E45").Value = '  +3.37%  '
$ws.Range("E46").Value = '  +1.74%  '
$ws.Range("E47").Value = '  +0.54%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.955'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.78%  '
$ws.Range("E49").Value = '  +3.02%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.049'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +15.90%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '111.32'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.33%  '
